$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6538461538461539
$ws.Range("C2").Value = 0.6938775510204082
$ws.Range("D2").Value = 0.6732673267326732
$ws.Range("E2").Value = 49

$ws.Range("B3").Value = 0.6341463414634146
$ws.Range("C3").Value = 0.5909090909090909
$ws.Range("D3").Value = 0.611764705882353
$ws.Range("E3").Value = 44

$ws.Range("B4").Value = 0.6451612903225806
$ws.Range("C4").Value = 0.6451612903225806
$ws.Range("D4").Value = 0.6451612903225806
$ws.Range("E4").Value = 0.6451612903225806

$ws.Range("B5").Value = 0.6439962476547842
$ws.Range("C5").Value = 0.6423933209647495
$ws.Range("D5").Value = 0.6425160163075131
$ws.Range("E5").Value = 93

$ws.Range("B6").Value = 0.6445258125037826
$ws.Range("C6").Value = 0.6451612903225806
$ws.Range("D6").Value = 0.6441693125669303
$ws.Range("E6").Value = 93

$ws.Range("B7").Value = 0.6481481481481481
$ws.Range("C7").Value = 0.7142857142857143
$ws.Range("D7").Value = 0.6796116504854369
$ws.Range("E7").Value = 49

$ws.Range("B8").Value = 0.6410256410256411
$ws.Range("C8").Value = 0.5681818181818182
$ws.Range("D8").Value = 0.6024096385542169
$ws.Range("E8").Value = 44

$ws.Range("B9").Value = 0.6451612903225806
$ws.Range("C9").Value = 0.6451612903225806
$ws.Range("D9").Value = 0.6451612903225806
$ws.Range("E9").Value = 0.6451612903225806

$ws.Range("B10").Value = 0.6445868945868947
$ws.Range("C10").Value = 0.6412337662337663
$ws.Range("D10").Value = 0.6410106445198269
$ws.Range("E10").Value = 93

$ws.Range("B11").Value = 0.6447783598321234
$ws.Range("C11").Value = 0.6451612903225806
$ws.Range("D11").Value = 0.6430859674212038
$ws.Range("E11").Value = 93

$ws.Range("B12").Value = 0.6440677966101694
$ws.Range("C12").Value = 0.7755102040816326
$ws.Range("D12").Value = 0.7037037037037036
$ws.Range("E12").Value = 49

$ws.Range("B13").Value = 0.6764705882352942
$ws.Range("C13").Value = 0.5227272727272727
$ws.Range("D13").Value = 0.5897435897435898
$ws.Range("E13").Value = 44

$ws.Range("B14").Value = 0.6559139784946236
$ws.Range("C14").Value = 0.6559139784946236
$ws.Range("D14").Value = 0.6559139784946236
$ws.Range("E14").Value = 0.6559139784946236

$ws.Range("B15").Value = 0.6602691924227317
$ws.Range("C15").Value = 0.6491187384044527
$ws.Range("D15").Value = 0.6467236467236467
$ws.Range("E15").Value = 93

$ws.Range("B16").Value = 0.6593981496371101
$ws.Range("C16").Value = 0.6559139784946236
$ws.Range("D16").Value = 0.6497870906473057
$ws.Range("E16").Value = 93

$ws.Range("B17").Value = 0.5609756097560976
$ws.Range("C17").Value = 0.9387755102040817
$ws.Range("D17").Value = 0.7022900763358778
$ws.Range("E17").Value = 49

$ws.Range("B18").Value = 0.7272727272727273
$ws.Range("C18").Value = 0.1818181818181818
$ws.Range("D18").Value = 0.2909090909090909
$ws.Range("E18").Value = 44

$ws.Range("B19").Value = 0.5806451612903226
$ws.Range("C19").Value = 0.5806451612903226
$ws.Range("D19").Value = 0.5806451612903226
$ws.Range("E19").Value = 0.5806451612903226

$ws.Range("B20").Value = 0.6441241685144125
$ws.Range("C20").Value = 0.5602968460111317
$ws.Range("D20").Value = 0.4965995836224844
$ws.Range("E20").Value = 93

$ws.Range("B21").Value = 0.6396538158929976
$ws.Range("C21").Value = 0.5806451612903226
$ws.Range("D21").Value = 0.5076582122629894
$ws.Range("E21").Value = 93

$ws.Range("B22").Value = 0.6226415094339622
$ws.Range("C22").Value = 0.673469387755102
$ws.Range("D22").Value = 0.6470588235294118
$ws.Range("E22").Value = 49

$ws.Range("B23").Value = 0.6
$ws.Range("C23").Value = 0.5454545454545454
$ws.Range("D23").Value = 0.5714285714285713
$ws.Range("E23").Value = 44

$ws.Range("B24").Value = 0.6129032258064516
$ws.Range("C24").Value = 0.6129032258064516
$ws.Range("D24").Value = 0.6129032258064516
$ws.Range("E24").Value = 0.6129032258064516

$ws.Range("B25").Value = 0.6113207547169811
$ws.Range("C25").Value = 0.6094619666048238
$ws.Range("D25").Value = 0.6092436974789915
$ws.Range("E25").Value = 93

$ws.Range("B26").Value = 0.6119293974437005
$ws.Range("C26").Value = 0.6129032258064516
$ws.Range("D26").Value = 0.6112767687720249
$ws.Range("E26").Value = 93
